$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 58056.223
$ws.Range("I62").Value = 64813.25
$ws.Range("K62").Value = 64813.25
$ws.Range("M62").Value = -64189.25

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 58056.223
$ws.Range("I65").Value = 64813.25
$ws.Range("K65").Value = 324066.25
$ws.Range("M65").Value = -320946.25

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1848.8043
$ws.Range("I132").Value = 1972.5897
$ws.Range("K132").Value = 5917.7691
$ws.Range("M132").Value = -3387.7691

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2185
$ws.Range("I137").Value = 1993.6666
$ws.Range("J137").Value = 3333
$ws.Range("K137").Value = 5980.9998
$ws.Range("L137").Value = 9999
$ws.Range("M137").Value = -3430.9998
$ws.Range("N137").Value = -15099

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3334.1892
$ws.Range("I138").Value = 2803.647
$ws.Range("J138").Value = 3785.15
$ws.Range("K138").Value = 8410.940999999999
$ws.Range("L138").Value = 11355.45
$ws.Range("M138").Value = -3270.940999999999
$ws.Range("N138").Value = -21635.45

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2951.1
$ws.Range("I2").Value = 2945.7778
$ws.Range("J2").Value = 2999
$ws.Range("K2").Value = 2945.7778
$ws.Range("L2").Value = 2999
$ws.Range("M2").Value = -2832.7778
$ws.Range("N2").Value = -3225

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4019.2222
$ws.Range("I61").Value = 1599.8572
$ws.Range("J61").Value = 5111.839
$ws.Range("K61").Value = 1599.8572
$ws.Range("L61").Value = 5111.839
$ws.Range("M61").Value = -1387.8572
$ws.Range("N61").Value = -5535.839

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2638.8
$ws.Range("I63").Value = 2638.8
$ws.Range("K63").Value = 2638.8
$ws.Range("M63").Value = -1952.8

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2638.8
$ws.Range("I66").Value = 2638.8
$ws.Range("K66").Value = 13194
$ws.Range("M66").Value = -9762

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1002424.6
$ws.Range("I74").Value = 2501727
$ws.Range("J74").Value = 2889.6667
$ws.Range("K74").Value = 2501727
$ws.Range("L74").Value = 2889.6667
$ws.Range("M74").Value = -2500853
$ws.Range("N74").Value = -4637.6667

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1002424.6
$ws.Range("I77").Value = 2501727
$ws.Range("J77").Value = 2889.6667
$ws.Range("K77").Value = 12508635
$ws.Range("L77").Value = 14448.3335
$ws.Range("M77").Value = -12504267
$ws.Range("N77").Value = -23184.3335

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2951.1
$ws.Range("I116").Value = 2945.7778
$ws.Range("J116").Value = 2999
$ws.Range("K116").Value = 2945.7778
$ws.Range("L116").Value = 2999
$ws.Range("M116").Value = -651.7777999999998
$ws.Range("N116").Value = -7587

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4019.2222
$ws.Range("I136").Value = 1599.8572
$ws.Range("J136").Value = 5111.839
$ws.Range("K136").Value = 4799.571599999999
$ws.Range("L136").Value = 15335.517
$ws.Range("M136").Value = -2249.571599999999
$ws.Range("N136").Value = -20435.517

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2951.1
$ws.Range("I3").Value = 2945.7778
$ws.Range("J3").Value = 2999
$ws.Range("K3").Value = 2945.7778
$ws.Range("L3").Value = 2999
$ws.Range("M3").Value = -2831.7778
$ws.Range("N3").Value = -3227

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1912.279
$ws.Range("I31").Value = 1899.9762
$ws.Range("K31").Value = 1899.9762
$ws.Range("M31").Value = -1604.9762

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1912.279
$ws.Range("I34").Value = 1899.9762
$ws.Range("K34").Value = 1899.9762
$ws.Range("M34").Value = -1697.9762

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 27052980
$ws.Range("I107").Value = 41704836
$ws.Range("J107").Value = 3399.2307
$ws.Range("K107").Value = 41704836
$ws.Range("L107").Value = 3399.2307
$ws.Range("M107").Value = -41702916
$ws.Range("N107").Value = -7239.2307

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7481.6
$ws.Range("I132").Value = 2380.423
$ws.Range("K132").Value = 7141.268999999999
$ws.Range("M132").Value = -4611.268999999999

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1638.738
$ws.Range("I134").Value = 1565.2646
$ws.Range("K134").Value = 4695.793799999999
$ws.Range("M134").Value = -2160.793799999999

# CUL row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 586.75
$ws.Range("I8").Value = 586.75
$ws.Range("K8").Value = 1760.25
$ws.Range("M8").Value = -1621.25

# CUL row 42
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 5000
$ws.Range("I42").Value = 5000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -14466
$ws.Range("N42").ClearContents()

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 859.2222
$ws.Range("J92").Value = 941.375
$ws.Range("L92").Value = 2824.125
$ws.Range("N92").Value = -5320.125

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3297.3333
$ws.Range("I131").Value = 1615.4445
$ws.Range("J131").Value = 5820.1665
$ws.Range("K131").Value = 4846.333500000001
$ws.Range("L131").Value = 17460.4995
$ws.Range("M131").Value = 193.6664999999994
$ws.Range("N131").Value = -27540.4995

# GSM row 95
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 175780.67
$ws.Range("J95").Value = 175780.67
$ws.Range("L95").Value = 175780.67
$ws.Range("N95").Value = -181272.67

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1916.7142
$ws.Range("I22").Value = 747.8333
$ws.Range("J22").Value = 2793.375
$ws.Range("K22").Value = 747.8333
$ws.Range("L22").Value = 2793.375
$ws.Range("M22").Value = -452.8333
$ws.Range("N22").Value = -3383.375

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1916.7142
$ws.Range("I27").Value = 747.8333
$ws.Range("J27").Value = 2793.375
$ws.Range("K27").Value = 747.8333
$ws.Range("L27").Value = 2793.375
$ws.Range("M27").Value = -640.8333
$ws.Range("N27").Value = -3007.375

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3010.2295
$ws.Range("I40").Value = 2207.2415
$ws.Range("K40").Value = 2207.2415
$ws.Range("M40").Value = -2071.2415

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4378.5
$ws.Range("I68").Value = 2299.6667
$ws.Range("K68").Value = 2299.6667
$ws.Range("M68").Value = -1550.6667

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4378.5
$ws.Range("I71").Value = 2299.6667
$ws.Range("K71").Value = 11498.3335
$ws.Range("M71").Value = -7754.333500000001

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1859.4667
$ws.Range("I82").Value = 1881
$ws.Range("J82").Value = 1816.4
$ws.Range("K82").Value = 1881
$ws.Range("L82").Value = 1816.4
$ws.Range("M82").Value = -1520
$ws.Range("N82").Value = -2538.4

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1859.4667
$ws.Range("I85").Value = 1881
$ws.Range("J85").Value = 1816.4
$ws.Range("K85").Value = 1881
$ws.Range("L85").Value = 1816.4
$ws.Range("M85").Value = -633
$ws.Range("N85").Value = -4312.4

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3749.4783
$ws.Range("I100").Value = 2520.611
$ws.Range("K100").Value = 2520.611
$ws.Range("M100").Value = -1979.611

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3378.653
$ws.Range("I122").Value = 3379.4783
$ws.Range("J122").Value = 3366
$ws.Range("K122").Value = 10138.4349
$ws.Range("L122").Value = 10098
$ws.Range("M122").Value = -7688.4349
$ws.Range("N122").Value = -14998

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2949.0513
$ws.Range("I132").Value = 2298.6272
$ws.Range("K132").Value = 6895.8816
$ws.Range("M132").Value = -4365.8816

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5529.9287
$ws.Range("I62").Value = 3752.1428
$ws.Range("K62").Value = 3752.1428
$ws.Range("M62").Value = -3128.1428

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5529.9287
$ws.Range("I65").Value = 3752.1428
$ws.Range("K65").Value = 18760.714
$ws.Range("M65").Value = -15640.714
